# Before the edit, the deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml -> clrScheme "Integral" (the theme the slide
#                            master / every slide actually designs against)
#   ppt/theme/theme2.xml -> clrScheme "Office"   (the stock default theme,
#                            only ever referenced by the notes master)
#
# The authored change swaps the contents of those two parts: the slide
# master's theme becomes the plain "Office Theme" palette, while the
# "Integral" palette that used to live there is relocated into the part the
# notes master points at. fontScheme/fmtScheme are identical between the two
# theme parts, so the entire observable change is the 12 clrScheme colors
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) plus the theme/clrScheme
# display names.
#
# Reproduce the part of that which is reachable through the exposed
# PowerPoint object model: ThemeColorScheme.Colors(i).RGB is the live COM
# surface that round-trips into the slide master's theme part
# (ppt/theme/theme1.xml), so drive the 12 slots to the standard "Office"
# theme RGB values. (There's no supported COM surface to rename the
# Design/ThemeColorScheme or to reach the notes-master-only theme part, so
# those calls are attempted defensively but don't block the color update if
# unsupported.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# slot -> new RGB, packed as R + G*256 + B*65536 (VBA RGB() order),
# matching the target "Office" clrScheme.
$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72

# Best-effort rename to match the target theme/clrScheme display names
# ("Integral" -> "Office Theme" / "Office"). Harmless if the host doesn't
# implement the setter.
try { $p.Designs.Item(1).Name = "Office Theme" } catch {}
try { $tcs.Name = "Office" } catch {}
